# v1.1: Improved visualisation readability
#
# The "MS" HealthCondition label is renamed to the fuller "Musculoskeletal"
# everywhere it appears in column C, the selection is moved onto the newly
# widened HealthCondition column, and that column is resized to fit its
# (now longer) contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename every "MS" HealthCondition entry (column C, including the data
# rows) to "Musculoskeletal".
$healthConditionRange = $ws.Range("C1:C33")
$healthConditionRange.Replace("MS", "Musculoskeletal")

# Resize column C (HealthCondition) to fit the new, longer text.
$ws.Columns.Item(3).AutoFit()

# Move the active selection onto the HealthCondition column range that was
# just updated.
$ws.Range("C18:C25").Select()
